# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#   - Status cells move from "Ready for handoff" to "Handed back: in sync with en-US"
#   - The Latest Handback DateTime for zh-cn / de-de is refreshed
#   - The stale "handback file is not the latest" error is cleared now that
#     the file is back in sync

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns (E2, F2)
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# zh-cn detail sheet
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("K2").Value = "2016-08-28 04:48:40"
$zhcn.Range("P2").Value = ""

# de-de detail sheet
$dede.Range("C2").Value = $newStatus
$dede.Range("K2").Value = "2016-08-28 04:48:47"
$dede.Range("P2").Value = ""

# Re-fit the widened Status / shrunk Error Detail columns like Excel would
# after the new, longer status text / now-empty error text is entered.
$overview.Columns.Item("E").AutoFit() | Out-Null
$overview.Columns.Item("F").AutoFit() | Out-Null
$zhcn.Columns.Item("C").AutoFit() | Out-Null
$zhcn.Columns.Item("P").AutoFit() | Out-Null
$dede.Columns.Item("C").AutoFit() | Out-Null
$dede.Columns.Item("P").AutoFit() | Out-Null

Write-Output "Handback report generated."
